$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

$years = 2013..2017
$startRow = 13

foreach ($ws in @($ws1, $ws2)) {
    $r = $startRow
    foreach ($year in $years) {
        $ws.Cells.Item($r, 1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\$year.xlsx"
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = "$year"
        $r++
    }
}

# Update the selected range on each sheet to cover the newly added rows,
# mirroring the original A2:B12 -> A2:B17 selection change. Hoja2 is done
# first so that Hoja1 ends up as the active/selected tab, matching the
# original workbook where Hoja1 has tabSelected="1".
$ws2.Activate()
$ws2.Range("A2:B17").Select()

$ws1.Activate()
$ws1.Range("A2:B17").Select()
